$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 5 ("B1- TYPES OF FINANCIAL DOCUMENTS"): the table (shape 2) gets a
#    different built-in table style applied (Table Design gallery pick).
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$tableShape = $s5.Shapes.Item(2)
$tbl = $tableShape.Table
$tbl.ApplyStyle("{BCC2B878-4213-410E-8983-233B2B9EE617}")

# ---------------------------------------------------------------------------
# 2) Theme swap: the deck's theme (the one backing the slide master / all
#    slides) changes from the "Integral" (Red Violet) palette to the plain
#    "Office Theme" palette. Re-point every theme color to the Office
#    palette value, in clrScheme order: dk1, lt1, dk2, lt2, accent1-6,
#    hlink, folHlink.
# ---------------------------------------------------------------------------
$tcs = $s5.ThemeColorScheme
$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
